# Auto-generated edit script: applies numeric cell updates to the
# Coeurl Profits workbook per the scheduled-runner refresh diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 376.8421
$ws.Range("J2").Value = 479.8
$ws.Range("L2").Value = 479.8
$ws.Range("N2").Value = -705.8
$ws.Range("H86").Value = 5614.6665
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 7672
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 7672
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -9918
$ws.Range("H89").Value = 5614.6665
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 7672
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 38360
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -49592
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992
$ws.Range("H103").Value = 530
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("H107").Value = 1339.4286
$ws.Range("I107").Value = 1398.3334
$ws.Range("J107").Value = 1295.25
$ws.Range("K107").Value = 1398.3334
$ws.Range("L107").Value = 1295.25
$ws.Range("M107").Value = 521.6666
$ws.Range("N107").Value = -5135.25
$ws.Range("H112").Value = 84885
$ws.Range("J112").Value = 145090.14
$ws.Range("L112").Value = 435270.42
$ws.Range("N112").Value = -437486.42
$ws.Range("H137").Value = 1558.6666
$ws.Range("I137").Value = 1470.5
$ws.Range("J137").Value = 1999.5
$ws.Range("K137").Value = 4411.5
$ws.Range("L137").Value = 5998.5
$ws.Range("M137").Value = -1861.5
$ws.Range("N137").Value = -11098.5
$ws.Range("H138").Value = 2893.09
$ws.Range("J138").Value = 3834.3052
$ws.Range("L138").Value = 11502.9156
$ws.Range("N138").Value = -21782.9156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5030.4546
$ws.Range("I32").Value = 4541.3
$ws.Range("K32").Value = 4541.3
$ws.Range("M32").Value = -4254.3
$ws.Range("H132").Value = 2475.1614
$ws.Range("I132").Value = 2073.8542
$ws.Range("J132").Value = 3851.0715
$ws.Range("K132").Value = 6221.562600000001
$ws.Range("L132").Value = 11553.2145
$ws.Range("M132").Value = -3691.562600000001
$ws.Range("N132").Value = -16613.2145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2070.721
$ws.Range("I105").Value = 1208.1613
$ws.Range("K105").Value = 1208.1613
$ws.Range("M105").Value = 538.8387
$ws.Range("H134").Value = 1432.9814
$ws.Range("I134").Value = 1488.1041
$ws.Range("K134").Value = 4464.3123
$ws.Range("M134").Value = -1929.3123

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4774.5
$ws.Range("I16").Value = 4033
$ws.Range("J16").Value = 6999
$ws.Range("K16").Value = 4033
$ws.Range("L16").Value = 6999
$ws.Range("M16").Value = -3746
$ws.Range("N16").Value = -7573
$ws.Range("H31").Value = 5526.381
$ws.Range("I31").Value = 2162.3635
$ws.Range("J31").Value = 9226.799999999999
$ws.Range("K31").Value = 2162.3635
$ws.Range("L31").Value = 9226.799999999999
$ws.Range("M31").Value = -1867.3635
$ws.Range("N31").Value = -9816.799999999999
$ws.Range("H34").Value = 5526.381
$ws.Range("I34").Value = 2162.3635
$ws.Range("J34").Value = 9226.799999999999
$ws.Range("K34").Value = 2162.3635
$ws.Range("L34").Value = 9226.799999999999
$ws.Range("M34").Value = -1960.3635
$ws.Range("N34").Value = -9630.799999999999
$ws.Range("H35").Value = 638.5454999999999
$ws.Range("I35").Value = 402.66666
$ws.Range("K35").Value = 402.66666
$ws.Range("M35").Value = -108.66666
$ws.Range("H58").Value = 2589.5134
$ws.Range("I58").Value = 2736.5356
$ws.Range("J58").Value = 2132.111
$ws.Range("K58").Value = 2736.5356
$ws.Range("L58").Value = 2132.111
$ws.Range("M58").Value = -2533.5356
$ws.Range("N58").Value = -2538.111
$ws.Range("H99").Value = 10708.429
$ws.Range("I99").Value = 9986.333000000001
$ws.Range("K99").Value = 9986.333000000001
$ws.Range("M99").Value = -8488.333000000001
$ws.Range("H113").Value = 4774.5
$ws.Range("I113").Value = 4033
$ws.Range("J113").Value = 6999
$ws.Range("K113").Value = 4033
$ws.Range("L113").Value = 6999
$ws.Range("M113").Value = -1863
$ws.Range("N113").Value = -11339
$ws.Range("H126").Value = 10708.429
$ws.Range("I126").Value = 9986.333000000001
$ws.Range("K126").Value = 29958.999
$ws.Range("M126").Value = -27488.999
$ws.Range("H132").Value = 3114.2766
$ws.Range("I132").Value = 2881.8462
$ws.Range("J132").Value = 4247.375
$ws.Range("K132").Value = 8645.5386
$ws.Range("L132").Value = 12742.125
$ws.Range("M132").Value = -6115.5386
$ws.Range("N132").Value = -17802.125
$ws.Range("H136").Value = 2589.5134
$ws.Range("I136").Value = 2736.5356
$ws.Range("J136").Value = 2132.111
$ws.Range("K136").Value = 8209.606800000001
$ws.Range("L136").Value = 6396.333
$ws.Range("M136").Value = -5659.606800000001
$ws.Range("N136").Value = -11496.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 556.0454999999999
$ws.Range("J113").Value = 623.5625
$ws.Range("L113").Value = 1870.6875
$ws.Range("N113").Value = -6210.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3634.3333
$ws.Range("I97").Value = 3201.5
$ws.Range("J97").Value = 4500
$ws.Range("K97").Value = 3201.5
$ws.Range("L97").Value = 4500
$ws.Range("M97").Value = -2705.5
$ws.Range("N97").Value = -5492
$ws.Range("H99").Value = 7192.5
$ws.Range("I99").Value = 7573.636
$ws.Range("K99").Value = 7573.636
$ws.Range("M99").Value = -5327.636
$ws.Range("H133").Value = 54500
$ws.Range("J133").Value = 54500
$ws.Range("L133").Value = 54500
$ws.Range("N133").Value = -64620

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1923.875
$ws.Range("I16").Value = 718.8
$ws.Range("J16").Value = 20000
$ws.Range("K16").Value = 718.8
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = -548.8
$ws.Range("N16").Value = -20340
$ws.Range("H46").Value = 2127.4
$ws.Range("J46").Value = 2823.75
$ws.Range("L46").Value = 2823.75
$ws.Range("N46").Value = -3199.75
$ws.Range("H68").Value = 2093.5
$ws.Range("I68").Value = 2127.9092
$ws.Range("K68").Value = 2127.9092
$ws.Range("M68").Value = -1378.9092
$ws.Range("H71").Value = 2093.5
$ws.Range("I71").Value = 2127.9092
$ws.Range("K71").Value = 10639.546
$ws.Range("M71").Value = -6895.546
$ws.Range("H95").Value = 34997.5
$ws.Range("J95").Value = 34997.5
$ws.Range("L95").Value = 34997.5
$ws.Range("N95").Value = -40489.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 100000
$ws.Range("I62").Value = 100000
$ws.Range("K62").Value = 100000
$ws.Range("M62").Value = -99376
$ws.Range("H65").Value = 100000
$ws.Range("I65").Value = 100000
$ws.Range("K65").Value = 500000
$ws.Range("M65").Value = -496880
$ws.Range("H107").Value = 407.9375
$ws.Range("I107").Value = 415.2
$ws.Range("J107").Value = 299
$ws.Range("K107").Value = 1245.6
$ws.Range("L107").Value = 897
$ws.Range("M107").Value = 674.4000000000001
$ws.Range("N107").Value = -4737
$ws.Range("H125").Value = 125025750
$ws.Range("J125").Value = 125025750
$ws.Range("L125").Value = 125025750
$ws.Range("N125").Value = -125035590
$ws.Range("H136").Value = 2260.5757
$ws.Range("I136").Value = 2058.6843
$ws.Range("K136").Value = 6176.0529
$ws.Range("M136").Value = -3626.0529
